$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5388.815
$ws.Range("I64").Value = 4722.6665
$ws.Range("J64").Value = 6721.1113
$ws.Range("K64").Value = 4722.6665
$ws.Range("L64").Value = 6721.1113
$ws.Range("M64").Value = -4474.6665
$ws.Range("N64").Value = -7217.1113
$ws.Range("H67").Value = 5388.815
$ws.Range("I67").Value = 4722.6665
$ws.Range("J67").Value = 6721.1113
$ws.Range("K67").Value = 4722.6665
$ws.Range("L67").Value = 6721.1113
$ws.Range("M67").Value = -3864.6665
$ws.Range("N67").Value = -8437.1113
$ws.Range("H69").Value = 1579.8182
$ws.Range("I69").Value = 1500
$ws.Range("J69").Value = 1939
$ws.Range("K69").Value = 4500
$ws.Range("L69").Value = 5817
$ws.Range("M69").Value = -3626
$ws.Range("N69").Value = -7565
$ws.Range("H72").Value = 1579.8182
$ws.Range("I72").Value = 1500
$ws.Range("J72").Value = 1939
$ws.Range("K72").Value = 13500
$ws.Range("L72").Value = 17451
$ws.Range("M72").Value = -9132
$ws.Range("N72").Value = -26187
$ws.Range("H131").Value = 2091.0715
$ws.Range("I131").Value = 2106.1538
$ws.Range("J131").Value = 1895
$ws.Range("K131").Value = 6318.4614
$ws.Range("L131").Value = 5685
$ws.Range("M131").Value = -1278.4614
$ws.Range("N131").Value = -15765
$ws.Range("H141").Value = 1886.9565
$ws.Range("I141").Value = 1766.8
$ws.Range("J141").Value = 2688
$ws.Range("K141").Value = 5300.4
$ws.Range("L141").Value = 8064
$ws.Range("M141").Value = -120.3999999999996
$ws.Range("N141").Value = -18424

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3048.8157
$ws.Range("I20").Value = 2721.3928
$ws.Range("J20").Value = 3965.6
$ws.Range("K20").Value = 2721.3928
$ws.Range("L20").Value = 3965.6
$ws.Range("M20").Value = -2474.3928
$ws.Range("N20").Value = -4459.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4006.3262
$ws.Range("I58").Value = 3376.4194
$ws.Range("J58").Value = 5308.1333
$ws.Range("K58").Value = 3376.4194
$ws.Range("L58").Value = 5308.1333
$ws.Range("M58").Value = -3173.4194
$ws.Range("N58").Value = -5714.1333
$ws.Range("H62").Value = 8255
$ws.Range("I62").Value = 4359.8
$ws.Range("J62").Value = 11501
$ws.Range("K62").Value = 4359.8
$ws.Range("L62").Value = 11501
$ws.Range("M62").Value = -3735.8
$ws.Range("N62").Value = -12749
$ws.Range("H65").Value = 8255
$ws.Range("I65").Value = 4359.8
$ws.Range("J65").Value = 11501
$ws.Range("K65").Value = 21799
$ws.Range("L65").Value = 57505
$ws.Range("M65").Value = -18679
$ws.Range("N65").Value = -63745
$ws.Range("H132").Value = 4049.5217
$ws.Range("I132").Value = 2886.2632
$ws.Range("J132").Value = 9575
$ws.Range("K132").Value = 8658.7896
$ws.Range("L132").Value = 28725
$ws.Range("M132").Value = -6128.7896
$ws.Range("N132").Value = -33785
$ws.Range("H136").Value = 4006.3262
$ws.Range("I136").Value = 3376.4194
$ws.Range("J136").Value = 5308.1333
$ws.Range("K136").Value = 10129.2582
$ws.Range("L136").Value = 15924.3999
$ws.Range("M136").Value = -7579.2582
$ws.Range("N136").Value = -21024.3999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 1549.8823
$ws.Range("I137").Value = 1039.6
$ws.Range("J137").Value = 1762.5
$ws.Range("K137").Value = 3118.8
$ws.Range("L137").Value = 5287.5
$ws.Range("M137").Value = 1981.2
$ws.Range("N137").Value = -15487.5
$ws.Range("H141").Value = 5434.84
$ws.Range("I141").Value = 5131.952
$ws.Range("J141").Value = 7025
$ws.Range("K141").Value = 15395.856
$ws.Range("L141").Value = 21075
$ws.Range("M141").Value = -10215.856
$ws.Range("N141").Value = -31435

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5146.3887
$ws.Range("I80").Value = 4538.6
$ws.Range("J80").Value = 5380.154
$ws.Range("K80").Value = 4538.6
$ws.Range("L80").Value = 5380.154
$ws.Range("M80").Value = -3540.6
$ws.Range("N80").Value = -7376.154
$ws.Range("H83").Value = 5146.3887
$ws.Range("I83").Value = 4538.6
$ws.Range("J83").Value = 5380.154
$ws.Range("K83").Value = 22693
$ws.Range("L83").Value = 26900.77
$ws.Range("M83").Value = -17701
$ws.Range("N83").Value = -36884.77
$ws.Range("H102").Value = 1816.6875
$ws.Range("I102").Value = 2296.5881
$ws.Range("J102").Value = 1553.5161
$ws.Range("K102").Value = 2296.5881
$ws.Range("L102").Value = 1553.5161
$ws.Range("M102").Value = -674.5880999999999
$ws.Range("N102").Value = -4797.5161

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 194.56522
$ws.Range("I55").Value = 97.86667
$ws.Range("J55").Value = 375.875
$ws.Range("K55").Value = 97.86667
$ws.Range("L55").Value = 375.875
$ws.Range("M55").Value = 75.13333
$ws.Range("N55").Value = -721.875
$ws.Range("H61").Value = 1288.4615
$ws.Range("I61").Value = 913.8
$ws.Range("J61").Value = 2537.3333
$ws.Range("K61").Value = 913.8
$ws.Range("L61").Value = 2537.3333
$ws.Range("M61").Value = -711.8
$ws.Range("N61").Value = -2941.3333
$ws.Range("H93").Value = 5098.8
$ws.Range("I93").Value = 4749
$ws.Range("J93").Value = 6498
$ws.Range("K93").Value = 4749
$ws.Range("L93").Value = 6498
$ws.Range("M93").Value = -3501
$ws.Range("N93").Value = -8994
$ws.Range("H113").Value = 1288.4615
$ws.Range("I113").Value = 913.8
$ws.Range("J113").Value = 2537.3333
$ws.Range("K113").Value = 913.8
$ws.Range("L113").Value = 2537.3333
$ws.Range("M113").Value = 1256.2
$ws.Range("N113").Value = -6877.3333
$ws.Range("H132").Value = 3664.8076
$ws.Range("I132").Value = 2738.5
$ws.Range("J132").Value = 4243.75
$ws.Range("K132").Value = 8215.5
$ws.Range("L132").Value = 12731.25
$ws.Range("M132").Value = -5685.5
$ws.Range("N132").Value = -17791.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8950.053
$ws.Range("I81").Value = 2653.2354
$ws.Range("J81").Value = 14047.477
$ws.Range("K81").Value = 5306.4708
$ws.Range("L81").Value = 28094.954
$ws.Range("M81").Value = -4245.4708
$ws.Range("N81").Value = -30216.954
$ws.Range("H84").Value = 8950.053
$ws.Range("I84").Value = 2653.2354
$ws.Range("J84").Value = 14047.477
$ws.Range("K84").Value = 26532.354
$ws.Range("L84").Value = 140474.77
$ws.Range("M84").Value = -21228.354
$ws.Range("N84").Value = -151082.77
$ws.Range("H113").Value = 1000.4211
$ws.Range("I113").Value = 682.3461
$ws.Range("J113").Value = 1689.5834
$ws.Range("K113").Value = 2047.0383
$ws.Range("L113").Value = 5068.7502
$ws.Range("M113").Value = 122.9617000000001
$ws.Range("N113").Value = -9408.7502
$ws.Range("H136").Value = 1997.6666
$ws.Range("I136").Value = 642.1142599999999
$ws.Range("J136").Value = 6742.1
$ws.Range("K136").Value = 1926.34278
$ws.Range("L136").Value = 20226.3
$ws.Range("M136").Value = 623.6572200000001
$ws.Range("N136").Value = -25326.3

